$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Intro paragraph: "For the assignment... within a bounding box. ... overall simulation box..."
Replace-Text "For the assignment, you will write code that simulates a set of particles moving randomly in a 2-dimensional space within a bounding box. The coordinates of the overall simulation box are between 0 and 100 along each dimension." `
             "For this assignment, you will write code to simulate a set of particles moving randomly in a 2-dimensional space within a global bounding box. The coordinates of the overall (global) simulation box are between 0 and 100 along each dimension."

# 2. "So, each chare owns a bounding box" -> "So, each chare owns a local bounding box"
Replace-Text "So, each chare owns a bounding box of its own with size" `
             "So, each chare owns a local bounding box of its own with size"

# 3. "within the chare) position" -> "within the chare's local bounding box) position"
Replace-Text "within the chare) position" `
             "within the chare’s local bounding box) position"

# 4. Reductions sentence
Replace-Text "every 10 iterations. The simulation should not be delayed by this calculation (you should use reductions)." `
             "every 10 iterations. The simulation should not be delayed by this calculation i.e. the simulation should not wait for this calculation and process (use asynchronous reductions)."

# 5. "use 10000 (=n) particles per chare, simulated over 100 steps(fixed) and a chare array"
Replace-Text "use 10000 (=n) particles per chare, simulated over 100 steps(fixed) and a chare array" `
             "use n=10000 particles per chare, simulated over 100 steps (fixed) and a chare array"

# 6. Move the "_GoBack" bookmark (Word tracks the last edit location) from the
# pseudocode paragraph to right before "(fixed)" in the sentence above, matching
# where the author's final edit in this revision landed.
$lastEdit = $d.Content
$lastEdit.Find.Execute("simulated over 100 steps ", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$goBackPoint = $d.Range($lastEdit.End, $lastEdit.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null
